$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-06 12:02:18", 0.0006000000000000001),
    @("2023-12-06 12:02:43", 0.0018),
    @("2023-12-06 12:03:39", 0.004),
    @("2023-12-06 12:03:46", 0.0004),
    @("2023-12-06 12:04:04", 0.0006000000000000001)
)

$row = 21
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row = $row + 1
}
